$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# DATE_TYPE_CODE (text "004" - write then reset style so no quote-prefix/number-format
# artifact is left behind, since the source keeps General formatting)
$ws.Range("J2").Value = "'004"
$ws.Range("J2").Style = "Normal"

# REPORT_DATE
$ws.Range("N2").Value = "2020-09-30 00:00:00"

# PARENT_NETPROFIT
$ws.Range("O2").Value = 958090978.66

# TOTAL_OPERATE_INCOME
$ws.Range("P2").Value = 15453781271.71

# TOTAL_OPERATE_COST
$ws.Range("Q2").Value = 14740290723.1

# TOE_RATIO
$ws.Range("R2").Value = -40.1675998012

# OPERATE_COST
$ws.Range("S2").Value = 11473644819.12

# OPERATE_EXPENSE
$ws.Range("T2").Value = 11473644819.12

# OPERATE_EXPENSE_RATIO
$ws.Range("U2").Value = -45.161077602

# SALE_EXPENSE
$ws.Range("V2").Value = 2392807895.23

# MANAGE_EXPENSE
$ws.Range("W2").Value = 739685498.8

# FINANCE_EXPENSE
$ws.Range("X2").Value = 24195678.38

# OPERATE_PROFIT
$ws.Range("Y2").Value = 1117370429.68

# TOTAL_PROFIT
$ws.Range("Z2").Value = 1125315017.31

# INCOME_TAX
$ws.Range("AA2").Value = 156025192.84

# OPERATE_TAX_ADD
$ws.Range("AG2").Value = 105920774.37

# TOI_RATIO
$ws.Range("AP2").Value = -39.57198276

# OPERATE_PROFIT_RATIO
$ws.Range("AQ2").Value = 5.816729740096

# PARENT_NETPROFIT_RATIO
$ws.Range("AR2").Value = 5.9

# DEDUCT_PARENT_NETPROFIT
$ws.Range("AS2").Value = 597717156.35

# DPN_RATIO
$ws.Range("AT2").Value = -31.162669673668
